$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force values to be written as literal TEXT,
# exactly preserving numeric-looking strings (e.g. "0.997", "69.953.99")
# instead of letting Excel auto-convert them to numbers.
$scratch = $ws.Range("Z1")

function Set-TextValue([string]$addr, [string]$text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

Set-TextValue 'D2' '69.953.99'
Set-TextValue 'E2' '  +0.79%  '
Set-TextValue 'D3' '3.536.47'
Set-TextValue 'E3' '  +1.25%  '
Set-TextValue 'D4' '0.997'
Set-TextValue 'E4' '  -0.13%  '
Set-TextValue 'D5' '606.49'
Set-TextValue 'E5' '  -1.08%  '
Set-TextValue 'D6' '197.38'
Set-TextValue 'E6' '  +5.98%  '
Set-TextValue 'D7' '0.631'
Set-TextValue 'E7' '  +0.34%  '
Set-TextValue 'D8' '0.998'
Set-TextValue 'E8' '  -0.11%  '
Set-TextValue 'E9' '  -6.71%  '
Set-TextValue 'D10' '0.651'
Set-TextValue 'E10' '  -0.12%  '
Set-TextValue 'D11' '53.91'
Set-TextValue 'E11' '  +1.11%  '
Set-TextValue 'D12' '0.0000303'
Set-TextValue 'E12' '  -1.43%  '
Set-TextValue 'D13' '9.54'
Set-TextValue 'E13' '  -0.47%  '
Set-TextValue 'D14' '4.085.33'
Set-TextValue 'E14' '  +0.87%  '
Set-TextValue 'D15' '600.45'
Set-TextValue 'E15' '  -0.91%  '
Set-TextValue 'D16' '69.984.89'
Set-TextValue 'E16' '  +0.84%  '
Set-TextValue 'D17' '19.12'
Set-TextValue 'E17' '  +1.46%  '
Set-TextValue 'D18' '12.72'
Set-TextValue 'E18' '  +0.60%  '
Set-TextValue 'D19' '3.515.94'
Set-TextValue 'E19' '  +0.46%  '
Set-TextValue 'D20' '0.122'
Set-TextValue 'E20' '  +1.32%  '
Set-TextValue 'D21' '0.995'
Set-TextValue 'E21' '  +0.61%  '
Set-TextValue 'D22' '18.49'
Set-TextValue 'E22' '  +7.04%  '
Set-TextValue 'D23' '5.29'
Set-TextValue 'E23' '  +4.85%  '
Set-TextValue 'E24' '  -3.23%  '
Set-TextValue 'D26' '3.20'
Set-TextValue 'E26' '  +5.74%  '
Set-TextValue 'D27' '10.92'
Set-TextValue 'E27' '  -0.30%  '
Set-TextValue 'D28' '9.62'
Set-TextValue 'E28' '  -3.60%  '
Set-TextValue 'E29' '  -0.72%  '
Set-TextValue 'D30' '4.37'
Set-TextValue 'E30' '  +11.47%  '
Set-TextValue 'D31' '7.11'
Set-TextValue 'E31' '  +1.52%  '
Set-TextValue 'D32' '12.53'
Set-TextValue 'E32' '  +0.54%  '
Set-TextValue 'E33' '  -1.91%  '
Set-TextValue 'D34' '63.09'
Set-TextValue 'E34' '  -0.36%  '
Set-TextValue 'D35' '0.0₃0859'
Set-TextValue 'E35' '  +10.80%  '
Set-TextValue 'D36' '3.725.94'
Set-TextValue 'E36' '  +4.15%  '
Set-TextValue 'E37' '  +0.13%  '
Set-TextValue 'E38' '  -3.05%  '
Set-TextValue 'E39' '  +1.40%  '
Set-TextValue 'D40' '0.393'
Set-TextValue 'D41' '36.68'
Set-TextValue 'E41' '  -0.25%  '
Set-TextValue 'D42' '488.24'
Set-TextValue 'E42' '  -6.61%  '
Set-TextValue 'D43' '0.133'
Set-TextValue 'E43' '  -3.48%  '
Set-TextValue 'D44' '0.0454'
Set-TextValue 'E44' '  -1.44%  '
Set-TextValue 'B45' 'Stellar'
Set-TextValue 'C45' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D45' '0.140'
Set-TextValue 'E45' '  -2.39%  '
Set-TextValue 'B46' 'ThetaToken'
Set-TextValue 'C46' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D46' '2.83'
Set-TextValue 'E46' '  -4.66%  '
Set-TextValue 'E47' '  -0.38%  '
Set-TextValue 'E48' '  +0.35%  '
Set-TextValue 'E49' '  -3.24%  '
Set-TextValue 'E50' '  +3.72%  '
Set-TextValue 'D51' '130.90'
Set-TextValue 'E51' '  -0.25%  '

$excel.CutCopyMode = $false
